$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# A new "% of total people" column is inserted before the existing
# "average debt per person" column, which now becomes column F.
# Copy the existing header formatting (bold / bordered / centered style)
# from E1 onto the new F1 header cell before changing E1's text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F1").Value = "average debt per person"
$ws.Range("E1").Value = "% of total people"

# --- Data rows ----------------------------------------------------------
# Occupations are now reported in descending order of headcount
# (merchant, gentleman, esq) and a "% of total people" figure has been
# added in column E, pushing the average-debt figure to column F.
$ws.Range("B2").Value = "merchant"
$ws.Range("C2").Value = 10577.739829
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1322.217478625

$ws.Range("B3").Value = "gentleman"
# C3 (4586.66) and D3 (1) already hold the correct figures for the
# "gentleman" row and do not need to change.
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 4586.66

$ws.Range("B4").Value = "esq"
$ws.Range("C4").Value = 102.089033
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 102.089033
